# Adds a small "FontFamily" / "FontSize" settings table to the
# "Einstellungen" sheet (columns M:N), matching the commit
# "Added font family and size to excel template file".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Einstellungen")
$ws.Activate()

# Fill the two header cells and the one data row for the new mini-table.
# Note: the header for column N ("FontSize") is written before the header
# for column M ("FontFamily") to match the original authoring order (this
# only affects shared-string ordering, not the final cell values).
$ws.Range("N1").Value = "FontSize"
$ws.Range("M1").Value = "FontFamily"
$ws.Range("M2").Value = "Arial, Helvetica, sans-serif"
$ws.Range("N2").Value = "0.6rem"

# Turn M1:N2 into a proper Excel Table ("Tabelle5"), headers already present.
$lo = $ws.ListObjects.Add(1, $ws.Range("M1:N2"), 0, 1)
$lo.Name = "Tabelle5"
$lo.TableStyle = "TableStyleLight12"

# Give the new column a bit more breathing room, like the other columns.
$ws.Columns.Item(13).ColumnWidth = 24

# Restore the view state close to what was captured in the saved file
# (scrolled a bit to the right, selection sitting on K21).
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("K21").Select()
